# v3.0 update FCI 27/1/2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header (13-01-2023 snapshot) - copy formatting from the existing
# B1 header cell so the new column matches the existing date-header style.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("C1").Value = "13-01-2023"

# Rows 2-8 are reordered (funds first, then avg, then total) and column C
# (13-01-2023 values) is populated alongside the existing column B values.
$rows = @(
    @{ Label = "Alpha Mega";           B = 47722.91;  C = 46505.29 },
    @{ Label = "Compass Small Cap II"; B = 57.76;     C = 46.83 },
    @{ Label = "Delta Acciones";       B = 10170.34;  C = 9739.07 },
    @{ Label = "Delta Select";         B = 302776.85; C = 303296.57 },
    @{ Label = "Delta gestion V";      B = 10709.79;  C = 8935.73 },
    @{ Label = "avg";                  B = 74287.53;  C = 73704.7 },
    @{ Label = "total";                B = 371437.65; C = 368523.49 }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.Label
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $r++
}
